# edit.ps1
# Applies the betexplorer turkey 1-lig 2023-2024 update:
#  1. Six pairs of adjacent rows had their match-detail columns (F:V) swapped
#     (the two fixtures played on the same date/time had their row order
#     reversed; columns A-E — index/country/tournament/season/date — stay put).
#  2. Three new fixture rows (140-142) were appended at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap columns F:V between each pair of rows that trade places.
# ---------------------------------------------------------------------------
$swapPairs = @(
    @(21, 22),
    @(31, 32),
    @(71, 72),
    @(96, 97),
    @(106, 107),
    @(121, 122)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $rangeA = $ws.Range("F$r1`:V$r1")
    $rangeB = $ws.Range("F$r2`:V$r2")
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

# ---------------------------------------------------------------------------
# 2) Append three new fixture rows (140, 141, 142) with the same formatting
#    (bold/bordered index column, date-formatted E column) as row 139.
# ---------------------------------------------------------------------------
$ws.Range("A139:V139").Copy()
$ws.Range("A140:V142").PasteSpecial(-4122)

$newRows = @(
    @{
        Row = 140; Idx = 139; Date = 45280.5
        Home = "Erzurumspor"; HomeGoals = 1; Away = "Umraniyespor"; AwayGoals = 1
        HomeOpenOdds = 2.8;  HomeOpenDt = "10/12/2023 11:43"; HomeCloseOdds = 2.17; HomeCloseDt = "20/12/2023 11:54"
        DrawOpenOdds = 3.11; DrawOpenDt = "10/12/2023 11:43"; DrawCloseOdds = 3.43; DrawCloseDt = "20/12/2023 11:54"
        AwayOpenOdds = 2.55; AwayOpenDt = "10/12/2023 11:43"; AwayCloseOdds = 3.41; AwayCloseDt = "20/12/2023 11:54"
        Url = "https://www.betexplorer.com/football/turkey/1-lig/erzurumspor-fk-umraniyespor/vmsp802I/"
    },
    @{
        Row = 141; Idx = 140; Date = 45280.5
        Home = "Keciorengucu"; HomeGoals = 1; Away = "Bandirmaspor"; AwayGoals = 3
        HomeOpenOdds = 3.02; HomeOpenDt = "10/12/2023 11:43"; HomeCloseOdds = 3.56; HomeCloseDt = "20/12/2023 11:56"
        DrawOpenOdds = 3.22; DrawOpenDt = "10/12/2023 11:43"; DrawCloseOdds = 3.43; DrawCloseDt = "20/12/2023 11:56"
        AwayOpenOdds = 2.31; AwayOpenDt = "10/12/2023 11:43"; AwayCloseOdds = 2.12; AwayCloseDt = "20/12/2023 11:55"
        Url = "https://www.betexplorer.com/football/turkey/1-lig/keciorengucu-bandirmaspor/lblXA2Ia/"
    },
    @{
        Row = 142; Idx = 141; Date = 45280.625
        Home = "Adanaspor AS"; HomeGoals = 0; Away = "Altay"; AwayGoals = 1
        HomeOpenOdds = 1.73; HomeOpenDt = "11/12/2023 18:13"; HomeCloseOdds = 1.66; HomeCloseDt = "20/12/2023 14:56"
        DrawOpenOdds = 3.61; DrawOpenDt = "11/12/2023 18:13"; DrawCloseOdds = 3.95; DrawCloseDt = "20/12/2023 14:59"
        AwayOpenOdds = 4.44; AwayOpenDt = "11/12/2023 18:13"; AwayCloseOdds = 5.12; AwayCloseDt = "20/12/2023 14:56"
        Url = "https://www.betexplorer.com/football/turkey/1-lig/adanaspor-as-altay/EudKD4Yt/"
    }
)

foreach ($rd in $newRows) {
    $r = $rd.Row
    $ws.Cells.Item($r, 1).Value2  = $rd.Idx          # A - Indice
    $ws.Cells.Item($r, 2).Value2  = "turkey"         # B - pais
    $ws.Cells.Item($r, 3).Value2  = "1-lig"          # C - torneio
    $ws.Cells.Item($r, 4).Value2  = "2023-2024"      # D - temporada
    $ws.Cells.Item($r, 5).Value2  = $rd.Date         # E - data_partida
    $ws.Cells.Item($r, 6).Value2  = $rd.Home         # F - home
    $ws.Cells.Item($r, 7).Value2  = $rd.HomeGoals    # G - home_ft_gols
    $ws.Cells.Item($r, 8).Value2  = $rd.Away         # H - away
    $ws.Cells.Item($r, 9).Value2  = $rd.AwayGoals    # I - away_ft_gols
    $ws.Cells.Item($r, 10).Value2 = $rd.HomeOpenOdds  # J - home_opening_odds
    $ws.Cells.Item($r, 11).Value2 = $rd.HomeOpenDt    # K - home_opening_data_hora
    $ws.Cells.Item($r, 12).Value2 = $rd.HomeCloseOdds # L - home_closing_odds
    $ws.Cells.Item($r, 13).Value2 = $rd.HomeCloseDt   # M - home_closing_data_hora
    $ws.Cells.Item($r, 14).Value2 = $rd.DrawOpenOdds  # N - draw_opening_odds
    $ws.Cells.Item($r, 15).Value2 = $rd.DrawOpenDt    # O - draw_opening_data_hora
    $ws.Cells.Item($r, 16).Value2 = $rd.DrawCloseOdds # P - draw_closing_odds
    $ws.Cells.Item($r, 17).Value2 = $rd.DrawCloseDt   # Q - draw_closing_data_hora
    $ws.Cells.Item($r, 18).Value2 = $rd.AwayOpenOdds  # R - away_opening_odds
    $ws.Cells.Item($r, 19).Value2 = $rd.AwayOpenDt    # S - away_opening_data_hora
    $ws.Cells.Item($r, 20).Value2 = $rd.AwayCloseOdds # T - away_closing_odds
    $ws.Cells.Item($r, 21).Value2 = $rd.AwayCloseDt   # U - away_closing_data_hora
    $ws.Cells.Item($r, 22).Value2 = $rd.Url           # V - url_partida
}

Write-Output "done"
